# Update the answer table in-place: each data row of the 5-column table
# gets its cell text replaced with the new values from the commit.
# Cell positions (row, col) are unchanged; only the <w:t> text differs.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (first data row)
$t.Cell(1,1).Range.Text = "77÷6=12, 5"
$t.Cell(1,2).Range.Text = "94÷6=15, 4"
$t.Cell(1,3).Range.Text = "44÷6=7, 2"
$t.Cell(1,4).Range.Text = "19÷6=3, 1"
$t.Cell(1,5).Range.Text = "33÷6=5, 3"

# Row 5 (second data row)
$t.Cell(5,1).Range.Text = "54÷3=18, 0"
$t.Cell(5,2).Range.Text = "53÷3=17, 2"
$t.Cell(5,3).Range.Text = "77÷9=8, 5"
$t.Cell(5,4).Range.Text = "91÷9=10, 1"
$t.Cell(5,5).Range.Text = "16÷4=4, 0"

# Row 9 (third data row)
$t.Cell(9,1).Range.Text = "44÷6=7, 2"
$t.Cell(9,2).Range.Text = "35÷8=4, 3"
$t.Cell(9,3).Range.Text = "19÷2=9, 1"
$t.Cell(9,4).Range.Text = "43÷9=4, 7"
$t.Cell(9,5).Range.Text = "86÷8=10, 6"

# Row 13 (fourth data row)
$t.Cell(13,1).Range.Text = "93÷5=18, 3"
$t.Cell(13,2).Range.Text = "52÷8=6, 4"
$t.Cell(13,3).Range.Text = "77÷7=11, 0"
$t.Cell(13,4).Range.Text = "22÷2=11, 0"
$t.Cell(13,5).Range.Text = "55÷5=11, 0"

# Row 17 (fifth data row)
$t.Cell(17,1).Range.Text = "50÷3=16, 2"
$t.Cell(17,2).Range.Text = "71÷8=8, 7"
$t.Cell(17,3).Range.Text = "62÷9=6, 8"
$t.Cell(17,4).Range.Text = "27÷4=6, 3"
$t.Cell(17,5).Range.Text = "76÷8=9, 4"
